$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells where the new value would otherwise be
# auto-converted to a number by Excel (losing formatting like trailing zeros,
# leading zeros, or precision).
$textCells = @("D5", "D6", "D10", "D11", "D15", "D16", "D20", "D21", "D23", "D25", "D29", "D31", "D32", "D33", "D34", "D35", "D37", "D39", "D41", "D42", "D43", "D45", "D46", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '62.971.35'
$ws.Range("E2").Value = '  +4.95%  '

$ws.Range("E3").Value = '  +5.69%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '573.46'
$ws.Range("E5").Value = '  +7.15%  '

$ws.Range("D6").Value = '152.76'
$ws.Range("E6").Value = '  +4.80%  '

$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").Value = '3.371.41'
$ws.Range("E8").Value = '  +5.54%  '

$ws.Range("E9").Value = '  -0.30%  '

$ws.Range("D10").Value = '7.44'
$ws.Range("E10").Value = '  +1.80%  '

$ws.Range("D11").Value = '0.119'
$ws.Range("E11").Value = '  +5.55%  '

$ws.Range("E12").Value = '  +1.08%  '

$ws.Range("D13").Value = '3.945.02'
$ws.Range("E13").Value = '  +5.66%  '

$ws.Range("E14").Value = '  +0.08%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '26.92'
$ws.Range("E15").Value = '  +4.00%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.0000180'
$ws.Range("E16").Value = '  +4.24%  '

$ws.Range("D17").Value = '63.044.33'
$ws.Range("E17").Value = '  +5.04%  '

$ws.Range("D18").Value = '3.382.21'
$ws.Range("E18").Value = '  +5.90%  '

$ws.Range("E19").Value = '  +0.59%  '

$ws.Range("D20").Value = '13.87'
$ws.Range("E20").Value = '  +4.89%  '

$ws.Range("D21").Value = '8.38'
$ws.Range("E21").Value = '  +2.09%  '

$ws.Range("E22").Value = '  +4.47%  '

$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("E24").Value = '  +2.37%  '

$ws.Range("D25").Value = '70.39'
$ws.Range("E25").Value = '  +1.36%  '

$ws.Range("E26").Value = '  +6.33%  '

$ws.Range("E27").Value = '  +6.91%  '

$ws.Range("D28").Value = '0.0₃0968'
$ws.Range("E28").Value = '  +11.32%  '

$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("E30").Value = '  +6.89%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '23.01'
$ws.Range("E31").Value = '  +3.01%  '

$ws.Range("B32").Value = 'RenderToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D32").Value = '6.35'
$ws.Range("E32").Value = '  +4.67%  '

$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = '5.58'
$ws.Range("E33").Value = '  +5.58%  '

$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = '1.31'
$ws.Range("E34").Value = '  +10.05%  '

$ws.Range("D35").Value = '6.70'
$ws.Range("E35").Value = '  +2.22%  '

$ws.Range("E36").Value = '  +9.45%  '

$ws.Range("D37").Value = '158.40'
$ws.Range("E37").Value = '  +1.72%  '

$ws.Range("E38").Value = '  +12.30%  '

$ws.Range("D39").Value = '27.36'
$ws.Range("E39").Value = '  +4.76%  '

$ws.Range("D40").Value = '2.889.07'
$ws.Range("E40").Value = '  +2.69%  '

$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '0.0740'
$ws.Range("E41").Value = '  +5.24%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.0327'
$ws.Range("E42").Value = '  +9.70%  '

$ws.Range("D43").Value = '40.80'
$ws.Range("E43").Value = '  +2.73%  '

$ws.Range("E44").Value = '  +4.43%  '

$ws.Range("D45").Value = '4.23'
$ws.Range("E45").Value = '  +0.29%  '

$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = '1.04'
$ws.Range("E46").Value = '  +5.62%  '

$ws.Range("B47").Value = 'RenzoRestakedETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D47").Value = '3.411.90'
$ws.Range("E47").Value = '  +5.71%  '

$ws.Range("D48").Value = '302.36'
$ws.Range("E48").Value = '  +14.55%  '

$ws.Range("D49").Value = '21.88'
$ws.Range("E49").Value = '  +5.58%  '

$ws.Range("E50").Value = '  -1.86%  '

$ws.Range("D51").Value = '6.28'
$ws.Range("E51").Value = '  +2.20%  '
